$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.382.94'
$ws.Range("E2").Value = '  +0.22%  '
$ws.Range("D3").Value = '3.505.23'
$ws.Range("E3").Value = '  -0.03%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.48'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.95%  '
$ws.Range("E10").Value = '  +0.06%  '
$ws.Range("E11").Value = '  +2.29%  '
$ws.Range("D12").Value = '4.105.41'
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("E13").Value = '  +1.64%  '
$ws.Range("E14").Value = '  +1.22%  '
$ws.Range("D15").Value = '3.508.04'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.80'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.69%  '
$ws.Range("D17").Value = '64.390.54'
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.92'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.76%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.76'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '393.33'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.575'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.33%  '
$ws.Range("D23").Value = '3.647.35'
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.62'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.07%  '
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.73'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000117'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.41'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.95%  '
$ws.Range("E30").Value = '  +2.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.28'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.47'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.159'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.79%  '
$ws.Range("D34").Value = '3.530.65'
$ws.Range("E34").Value = '  +0.28%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.42'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.67%  '
$ws.Range("E37").Value = '  +1.68%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.94'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.27%  '
$ws.Range("E39").Value = '  +0.81%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '167.53'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.73%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0787'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.82%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.810'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.14%  '
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '25.15'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.77%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.44'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.26%  '
$ws.Range("E46").Value = '  +5.69%  '
$ws.Range("E47").Value = '  -3.41%  '
$ws.Range("E48").Value = '  +0.78%  '
$ws.Range("D49").Value = '2.391.34'
$ws.Range("E49").Value = '  -3.53%  '
$ws.Range("E50").Value = '  -0.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0259'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.23%  '
